$wb = $excel.ActiveWorkbook

# Add the new worksheet "doFindStoreLocator" at the end
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "doFindStoreLocator"

$ws2.Range("A1").Value = "city"
$ws2.Range("A2").Value = "Bridgeport, CT"
$ws2.Range("A3").Value = "Lake Charles, LA"
$ws2.Range("A4").Value = "Brooklyn, NY"
$ws2.Range("A5").Value = "Abilene, TX"

$ws2.Columns.Item(1).ColumnWidth = 14

$ws2.Activate()
